# Update "想去人数" (want-to-go count) figures in column F on both the
# "展览" sheet and the mirrored "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 121
$ws1.Range("F14").Value = 185
$ws1.Range("F15").Value = 985
$ws1.Range("F18").Value = 154
$ws1.Range("F23").Value = 5879
$ws1.Range("F27").Value = 523
$ws1.Range("F29").Value = 3360
$ws1.Range("F41").Value = 910
$ws1.Range("F43").Value = 24

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 121
$ws4.Range("F15").Value = 185
$ws4.Range("F16").Value = 985
$ws4.Range("F19").Value = 154
$ws4.Range("F24").Value = 5879
$ws4.Range("F28").Value = 523
$ws4.Range("F30").Value = 3360
$ws4.Range("F42").Value = 910
$ws4.Range("F44").Value = 24
